$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A75").Value = "Globo"
$ws.Range("B75").Value = "Bom Dia Inter"
$ws.Range("C75").Value = "Obras"
$ws.Range("D75").Value = "2025-04-08T08:10"
$ws.Range("E75").Value = "Negativo"
$ws.Range("F75").Value = "Cuidado com os buracos. Em Campos, buracos no Parque Esplanada e Guarus causam riscos. Repórter *ao vivo* da Rua Gregório Marcolino, no Esplanada. Duas crateras se abriram. Segundo moradores, foram feitos reparos mas, depois, o buraco voltou a abrir. Aos poucos, está cedendo mais. Em Guarus, o problema é na Rua Hipólito Sardinha, perto da delegacia de Guarus. Tampa do bueiro afundou. Problema ocorre há cerca de dois meses, segundo entrevista dos moradores. Motoristas reclamam do descaso das autoridades. Em nota, Prefeitura informou que a Secretaria de Obras enviou técnico à rua de Guarus para avaliar a situação. *nota da Secretaria de Obras*"
